$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 19.67207418041587

# Row 3
$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 3.005019366241741

# Row 4
$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 5.582307763322248

# Row 5
$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 9.983522426115931
$ws.Range("D5").Value = 189.6080260415259
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 216.727722176622
